$d = $word.ActiveDocument

# wdAlignParagraphJustify = 3
$wdAlignParagraphJustify = 3

# --- Paragraph 1: 'Dear Editors and Reviewers,' -> add justify alignment ---
$d.Paragraphs.Item(1).Range.ParagraphFormat.Alignment = $wdAlignParagraphJustify

# --- Paragraph 2: empty paragraph -> add justify alignment ---
$d.Paragraphs.Item(2).Range.ParagraphFormat.Alignment = $wdAlignParagraphJustify

# --- Paragraph 3: 'We would like to thank you...' -> add justify + replace text ---
$d.Paragraphs.Item(3).Range.ParagraphFormat.Alignment = $wdAlignParagraphJustify
$d.Content.Find.Execute('We would like to thank you for the unanimous interest you showed for our study. The theoretical and analytical advice and criticism have helped us reshape the manuscript in a way which we hope will make it more conceptually and biologically relevant.', $true, $false, $false, $false, $false, $true, 1, $false, 'We would like to thank you for the unanimous interest you showed for our study. The theoretical and analytical advice and criticism have helped us reanalyze the data, and reshape the manuscript in a way which we hope will make it more conceptually and biologically relevant. We have tried here to synthesize the large number of changes made to the manuscript and their relation to comments and criticisms.', 2) | Out-Null

# --- Paragraph 4: 'Important changes...' -> add justify + replace text ---
$d.Paragraphs.Item(4).Range.ParagraphFormat.Alignment = $wdAlignParagraphJustify
$d.Content.Find.Execute('Important changes have been made in the Introduction, in which we added more details on our reasoning behind the idea of optimal asymmetry (lines 100-120), as asked by Reviewer 2. In addition, we have reworded the following paragraphs, including hypotheses and expectations (lines 122-171), so as to (hopefully) clarify and rectify our meaning behind modularity and integration, and how the covariation patterns observed may relate to those.', $true, $false, $false, $false, $false, $true, 1, $false, 'Important changes have been made in the Introduction, especially the second half, in which we added more details on our reasoning behind the idea of optimal asymmetry, as asked by Reviewer 2. In addition, we have reworked the following paragraphs, including hypotheses and expectations, so as to clarify and rectify our meaning behind modularity and integration, and how the covariation patterns observed may relate to those. We have kept on using the terms variational / developmental / functional modularity, because we feel that, although modularity in biology has been used mostly in the realm of developmental studies, it actually is a more general property of systems which we think applies in several ways in biology. The paragraph in the Methods section which was also referring to the different hypotheses, but was deemed unclear by Dr. Zelditch, has been removed. We hope that gathering all of these hypotheses and expectations in one place at the end of the Introduction, with better delineated paragraphs will make it easier for readers to follow our reasoning.', 2) | Out-Null

# --- Paragraph 5: 'Regarding the analyses...' -> add justify + replace text ---
$d.Paragraphs.Item(5).Range.ParagraphFormat.Alignment = $wdAlignParagraphJustify
$d.Content.Find.Execute('Regarding the analyses and results, we opted to remove the sections using “module by module” superimposition, as suggested by Dr. Collyer and Dr. Zelditch. While we did not want to appear dismissive of this recent literature, we also agree that these additional analyses somewhat disrupted the flow of the manuscript.', $true, $false, $false, $false, $false, $true, 1, $false, 'Regarding the analyses and results, we opted to remove entirely the sections using “module by module” superimposition, as suggested by Dr. Collyer and Dr. Zelditch, as well as the EMMLi analyses, which were indeed already relegated to the supplementary material. While we did not want to appear dismissive of this recent literature, we also agree that these additional analyses somewhat disrupted the flow of the reading, and we feel that their removal makes the manuscript more straightforward and to the point. As a consequence of this removal, some of the related remarks by Dr. Zelditch do not apply, because we do not discuss the differences between global and module-by-module superimpositions anymore.', 2) | Out-Null

# --- Append new paragraphs after paragraph 5, each justified ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $lastPara.Range
# new paragraph (source index 5)
$insertionPoint.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$insertionPoint.Text = 'Following Dr. Collyer’s reasoning, and the remarks of Reviewer 2 on the inadequacy of using modularity as a way to compare between partitioning hypotheses, we ran analyses with the head capsule and mandibles as separate structures. Instead of testing and trying to compare different partitions, we opted to start from the assumption that both mandibles and the dorsal and ventral parts of the head constitute modules, which is justified at least anatomically and developmentally, and to look at patterns of covariation as evidence supporting this or not, or suggesting that other processes may produce observed patterns. We have used covariance ratio and 2B-PLS analyses viewing them purely as tests for covariation between blocks and within blocks, but refrained from calling them “modularity” or “integration” tests. In addition, as suggested by Dr. Collyer, we computed and plotted the ordered covariance matrix for a more exhaustive and qualitative interpretation of covariation patterns within and between the proposed modules. In addition to the covariance matrix, we computed and plotted the congruence matrix, a matrix of 3D vector correlations based on the congruence coefficient. Both matrices are similar, but the latter has a readability advantage, because it is a landmark-by-landmark rather than coordinate-by-coordinate matrix, which in our opinion makes interpretation of covariation a bit easier for the reader. Both matrices were plotted as heatmaps, and show similar overall results. They were both represented in the new Figure 2.'
$insertionPoint = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$insertionPoint.ParagraphFormat.Alignment = $wdAlignParagraphJustify

# new paragraph (source index 6)
$insertionPoint.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$insertionPoint.Text = 'As suggested by Dr. Collyer, the relationships between mandible shapes and bite force were assessed via 2B-PLS. Furthermore, as suggested by reviewer 2, asymmetric shape differences are now represented as 3D morphs in the new Figure 3. Although this will be certainly helpful for the readers unfamiliar with insect head morphology, we felt that lollipop plots allowed for more details to be seen, of interest perhaps to the more specialist audience, we therefore opted to keep these as supplementary figures.'
$insertionPoint = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$insertionPoint.ParagraphFormat.Alignment = $wdAlignParagraphJustify

# new paragraph (source index 7)
$insertionPoint.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$insertionPoint.Text = 'In accordance with new results and with bibliographical suggestions made notably by reviewer 2, we have also largely modified the Discussion. We now orient the discussion on the various processes which may be at play to produce the observed covariation patterns, and propose that patterns related to developmental modularity may be obscured notably by variation due to plastic changes or wear related to the feeding function. Following this reasoning, we introduce the palimpsest model proposed by Hallgrímsson et al. (2009) as a conceptual framework to interpret our results. In addition, we have added discussion on the potential importance of selection for integration or parcellation in the evolution of bilateral modularity and asymmetry.'
$insertionPoint = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$insertionPoint.ParagraphFormat.Alignment = $wdAlignParagraphJustify

# new paragraph (source index 8)
$insertionPoint.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$insertionPoint.Text = 'Finally, the title and abstract have been modified to better match the new results and discussion of the study.'
$insertionPoint = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$insertionPoint.ParagraphFormat.Alignment = $wdAlignParagraphJustify

# new paragraph (source index 9)
$insertionPoint.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$insertionPoint.ParagraphFormat.Alignment = $wdAlignParagraphJustify

# new paragraph (source index 10)
$insertionPoint.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$insertionPoint.Text = 'We thank once more the reviewers and editors for their comments, and hope that the changes made to the manuscript adequately address them. We hope our modified manuscript will be deemed adequate for further consideration for publication in Evolution.'
$insertionPoint = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$insertionPoint.ParagraphFormat.Alignment = $wdAlignParagraphJustify

# new paragraph (source index 11)
$insertionPoint.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$insertionPoint.Text = 'Best regards'
$insertionPoint = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$insertionPoint.ParagraphFormat.Alignment = $wdAlignParagraphJustify

